$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) target cells to Text format to preserve exact string values
$priceCells = @('D2', 'D3', 'D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D12', 'D13', 'D14', 'D15', 'D16', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D29', 'D30', 'D31', 'D32', 'D33', 'D35', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D50', 'D51')
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range('D2').Value = '64.866.80'
$ws.Range('E2').Value = '  +0.79%  '
$ws.Range('D3').Value = '3.377.53'
$ws.Range('E3').Value = '  -0.95%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '559.25'
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('D6').Value = '177.34'
$ws.Range('E6').Value = '  +2.30%  '
$ws.Range('D7').Value = '0.623'
$ws.Range('E7').Value = '  +0.49%  '
$ws.Range('D8').Value = '3.369.21'
$ws.Range('E8').Value = '  -0.96%  '
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').Value = '0.168'
$ws.Range('E10').Value = '  +8.57%  '
$ws.Range('E11').Value = '  +2.05%  '
$ws.Range('D12').Value = '55.23'
$ws.Range('E12').Value = '  -2.91%  '
$ws.Range('D13').Value = '0.0000279'
$ws.Range('E13').Value = '  +3.44%  '
$ws.Range('D14').Value = '9.16'
$ws.Range('E14').Value = '  +1.32%  '
$ws.Range('D15').Value = '3.906.64'
$ws.Range('E15').Value = '  -1.31%  '
$ws.Range('D16').Value = '18.37'
$ws.Range('E16').Value = '  +2.09%  '
$ws.Range('E17').Value = '  -1.74%  '
$ws.Range('D18').Value = '3.369.92'
$ws.Range('E18').Value = '  -1.38%  '
$ws.Range('D19').Value = '11.92'
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('D20').Value = '64.715.76'
$ws.Range('E20').Value = '  +0.54%  '
$ws.Range('D21').Value = '0.991'
$ws.Range('E21').Value = '  +0.42%  '
$ws.Range('D22').Value = '459.64'
$ws.Range('E22').Value = '  +12.59%  '
$ws.Range('D23').Value = '4.78'
$ws.Range('E23').Value = '  +11.53%  '
$ws.Range('D24').Value = '4.11'
$ws.Range('E24').Value = '  -0.53%  '
$ws.Range('D25').Value = '86.05'
$ws.Range('E25').Value = '  +3.68%  '
$ws.Range('D26').Value = '13.51'
$ws.Range('E26').Value = '  +1.07%  '
$ws.Range('D27').Value = '10.94'
$ws.Range('E27').Value = '  +1.89%  '
$ws.Range('E28').Value = '  +3.54%  '
$ws.Range('D29').Value = '8.84'
$ws.Range('E29').Value = '  -0.30%  '
$ws.Range('D30').Value = '30.18'
$ws.Range('E30').Value = '  +1.96%  '
$ws.Range('D31').Value = '6.78'
$ws.Range('E31').Value = '  +1.78%  '
$ws.Range('D32').Value = '11.53'
$ws.Range('E32').Value = '  +0.64%  '
$ws.Range('D33').Value = '583.40'
$ws.Range('E33').Value = '  -1.24%  '
$ws.Range('E34').Value = '  +1.40%  '
$ws.Range('D35').Value = '59.59'
$ws.Range('E35').Value = '  +1.25%  '
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('E37').Value = '  -6.84%  '
$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').Value = '0.0₃0770'
$ws.Range('E38').Value = '  +3.82%  '
$ws.Range('B39').Value = 'InjectiveProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D39').Value = '36.07'
$ws.Range('E39').Value = '  +0.97%  '
$ws.Range('D40').Value = '3.47'
$ws.Range('E40').Value = '  +1.37%  '
$ws.Range('D41').Value = '0.374'
$ws.Range('E41').Value = '  +0.72%  '
$ws.Range('D42').Value = '3.106.43'
$ws.Range('E42').Value = '  -2.06%  '
$ws.Range('D43').Value = '0.999'
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').Value = '2.86'
$ws.Range('E44').Value = '  -0.96%  '
$ws.Range('D45').Value = '2.53'
$ws.Range('E45').Value = '  +0.18%  '
$ws.Range('D46').Value = '0.0414'
$ws.Range('E46').Value = '  +1.80%  '
$ws.Range('D47').Value = '3.22'
$ws.Range('E47').Value = '  -0.18%  '
$ws.Range('D48').Value = '0.132'
$ws.Range('E48').Value = '  +1.81%  '
$ws.Range('E49').Value = '  -1.81%  '
$ws.Range('D50').Value = '8.34'
$ws.Range('E50').Value = '  +0.79%  '
$ws.Range('D51').Value = '135.83'
$ws.Range('E51').Value = '  +0.65%  '

# Restore default (Normal) style on the price cells so no stray number-format styling remains
foreach ($cell in $priceCells) {
    $ws.Range($cell).Style = "Normal"
}
